# Update the cryptos list with latest scraped values (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a numeric-looking string as literal text,
    # matching the source data which stores prices as inline strings.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "38.027.86"
$ws.Range("E2").Value = "  +0.21%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.040.54"
$ws.Range("E3").Value = "  -0.58%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "228.90"
$ws.Range("E5").Value = "  -0.40%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.613"
$ws.Range("E6").Value = "  -0.69%  "

# Row 7 - Solana
Set-TextValue $ws.Range("D7") "60.68"
$ws.Range("E7").Value = "  +3.57%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.382"
$ws.Range("E9").Value = "  -0.78%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.19%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.92%  "

# Row 12 - Chainlink
Set-TextValue $ws.Range("D12") "14.66"
$ws.Range("E12").Value = "  +0.26%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.340.16"
$ws.Range("E13").Value = "  -0.65%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "21.41"
$ws.Range("E14").Value = "  +2.86%  "

# Row 15 - Polygon
Set-TextValue $ws.Range("D15") "0.767"
$ws.Range("E15").Value = "  +1.91%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  -1.69%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.035.96"
$ws.Range("E17").Value = "  -1.23%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "37.963.92"
$ws.Range("E18").Value = "  +0.27%  "

# Row 19/20 - Uniswap and Litecoin swap order (Litecoin now ranked 19, Uniswap 20)
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D19") "69.92"
$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D20") "5.97"
$ws.Range("E20").Value = "  -5.39%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  -1.08%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "224.83"
$ws.Range("E22").Value = "  +0.17%  "

# Row 23 - Dai
Set-TextValue $ws.Range("D23") "0.999"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.01%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -0.29%  "

# Row 26 - Cosmos
Set-TextValue $ws.Range("D26") "9.37"
$ws.Range("E26").Value = "  +0.94%  "

# Row 27 - Monero
Set-TextValue $ws.Range("D27") "167.38"
$ws.Range("E27").Value = "  +0.55%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -1.97%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "18.93"
$ws.Range("E29").Value = "  -0.51%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  -3.43%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  +0.73%  "

# Row 32 - WEMIXToken
$ws.Range("E32").Value = "  +5.53%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -2.17%  "

# Row 34/35 - Hedera and InternetComputer(DFINITY) swap order (ICP now ranked 34, Hedera 35)
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D34") "4.54"
$ws.Range("E34").Value = "  -1.16%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D35") "0.0608"
$ws.Range("E35").Value = "  -0.08%  "

# Row 36 - THORChain
Set-TextValue $ws.Range("D36") "6.50"
$ws.Range("E36").Value = "  +7.11%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -1.84%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  +0.45%  "

# Row 39 - BinanceUSD
$ws.Range("E39").Value = "  -0.23%  "

# Row 40 - InjectiveProtocol
Set-TextValue $ws.Range("D40") "17.78"
$ws.Range("E40").Value = "  +7.12%  "

# Row 41 - Maker
$ws.Range("D41").Value = "1.527.76"
$ws.Range("E41").Value = "  +1.81%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  +0.62%  "

# Row 43 - Aave
Set-TextValue $ws.Range("D43") "96.34"
$ws.Range("E43").Value = "  -0.79%  "

# Row 44 - HuobiToken
Set-TextValue $ws.Range("D44") "2.82"
$ws.Range("E44").Value = "  -2.33%  "

# Row 45 - Cronos
Set-TextValue $ws.Range("D45") "0.0917"
$ws.Range("E45").Value = "  -0.07%  "

# Row 46 - TrustWalletToken
$ws.Range("E46").Value = "  -2.78%  "

# Row 47 - FTXToken
Set-TextValue $ws.Range("D47") "4.05"
$ws.Range("E47").Value = "  -2.05%  "

# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  -0.33%  "

# Row 49 - MXToken
$ws.Range("E49").Value = "  +0.21%  "

# Row 50 - FraxShare
Set-TextValue $ws.Range("D50") "7.12"
$ws.Range("E50").Value = "  +0.30%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.230.38"
$ws.Range("E51").Value = "  -0.58%  "
